$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds line1..line6 (rows 2-7) followed by extr1..extr8
# (rows 8-15). Two new lines ("line7", "line8") need to be inserted right
# after line6, which pushes extr1..extr8 down to rows 10-17.
#
# Rather than using Rows.Insert() (which drags along border/format side
# effects on the newly created rows), shift the existing extr1..extr8 block
# down by writing its current values into rows 10-17, working from the
# bottom up so the source cells are read before they get overwritten.
for ($i = 7; $i -ge 0; $i--) {
    $src = 8 + $i
    $dst = 10 + $i
    $ws.Cells.Item($dst, 1).Value = $ws.Cells.Item($src, 1).Value2
    $ws.Cells.Item($dst, 2).Value = $ws.Cells.Item($src, 2).Value2
    $ws.Cells.Item($dst, 3).Value = $ws.Cells.Item($src, 3).Value2
    $ws.Cells.Item($dst, 4).Value = $ws.Cells.Item($src, 4).Value2
    $ws.Cells.Item($dst, 5).Value = $ws.Cells.Item($src, 5).Value2
}

# Give the two brand-new rows (16-17) the same formatting as the rest of the
# index column (bold, centered, bordered) by copying an existing cell's
# format onto them.
$ws.Range("A2").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Now fill in the two new rows with the line7 / line8 data (row 8 and 9
# already carry the correct index-column style from the original sheet).
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# Restore the running 0-based index in column A for every shifted data row
# (A = row - 2).
for ($r = 10; $r -le 17; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# Two in_service flags flip from FALSE to TRUE on the shifted extr4/extr5
# rows (now rows 13 and 14).
$ws.Range("E13").Value = $true
$ws.Range("E14").Value = $true
